$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.InsertAfter("Removed Log.java, Profile.java, and Invitation.java`r")
